# Insert 4 new weekly price rows above row 936 (shifts existing rows 936:1037
# down to 940:1041) and populate the newly inserted rows with the new week's
# data (Fecha = 45194), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 936 - everything currently at/after row 936 shifts
# down by 4 rows (old 936 -> new 940, ..., old 1037 -> new 1041).
$ws.Rows("936:939").Insert()

# Row 936: Lane Late / Primera
$ws.Range("A936").Value = 8
$ws.Range("B936").Value = "Terminal La Palmera de La Serena"
$ws.Range("C936").Value = "Coquimbo"
$ws.Range("D936").Value = 45194
$ws.Range("E936").Value = 4
$ws.Range("F936").Value = "Fruta"
$ws.Range("G936").Value = 100102
$ws.Range("H936").Value = "Cítricos"
$ws.Range("I936").Value = 100102005
$ws.Range("J936").Value = "Naranja"
$ws.Range("K936").Value = "Lane Late"
$ws.Range("L936").Value = "Primera"
$ws.Range("M936").Value = 20
$ws.Range("N936").Value = 140000
$ws.Range("O936").Value = 150000
$ws.Range("P936").Value = 145000
$ws.Range("Q936").Value = "$/bins (400 kilos)"
$ws.Range("R936").Value = "Provincia de Limarí"
$ws.Range("S936").Value = 362
$ws.Range("T936").Value = 400

# Row 937: Lane Late / Segunda
$ws.Range("A937").Value = 8
$ws.Range("B937").Value = "Terminal La Palmera de La Serena"
$ws.Range("C937").Value = "Coquimbo"
$ws.Range("D937").Value = 45194
$ws.Range("E937").Value = 4
$ws.Range("F937").Value = "Fruta"
$ws.Range("G937").Value = 100102
$ws.Range("H937").Value = "Cítricos"
$ws.Range("I937").Value = 100102005
$ws.Range("J937").Value = "Naranja"
$ws.Range("K937").Value = "Lane Late"
$ws.Range("L937").Value = "Segunda"
$ws.Range("M937").Value = 20
$ws.Range("N937").Value = 120000
$ws.Range("O937").Value = 130000
$ws.Range("P937").Value = 125000
$ws.Range("Q937").Value = "$/bins (400 kilos)"
$ws.Range("R937").Value = "Provincia de Limarí"
$ws.Range("S937").Value = 312
$ws.Range("T937").Value = 400

# Row 938: Navel Late / Primera
$ws.Range("A938").Value = 8
$ws.Range("B938").Value = "Terminal La Palmera de La Serena"
$ws.Range("C938").Value = "Coquimbo"
$ws.Range("D938").Value = 45194
$ws.Range("E938").Value = 4
$ws.Range("F938").Value = "Fruta"
$ws.Range("G938").Value = 100102
$ws.Range("H938").Value = "Cítricos"
$ws.Range("I938").Value = 100102005
$ws.Range("J938").Value = "Naranja"
$ws.Range("K938").Value = "Navel Late"
$ws.Range("L938").Value = "Primera"
$ws.Range("M938").Value = 20
$ws.Range("N938").Value = 140000
$ws.Range("O938").Value = 150000
$ws.Range("P938").Value = 145000
$ws.Range("Q938").Value = "$/bins (400 kilos)"
$ws.Range("R938").Value = "Provincia de Limarí"
$ws.Range("S938").Value = 362
$ws.Range("T938").Value = 400

# Row 939: Navel Late / Segunda
$ws.Range("A939").Value = 8
$ws.Range("B939").Value = "Terminal La Palmera de La Serena"
$ws.Range("C939").Value = "Coquimbo"
$ws.Range("D939").Value = 45194
$ws.Range("E939").Value = 4
$ws.Range("F939").Value = "Fruta"
$ws.Range("G939").Value = 100102
$ws.Range("H939").Value = "Cítricos"
$ws.Range("I939").Value = 100102005
$ws.Range("J939").Value = "Naranja"
$ws.Range("K939").Value = "Navel Late"
$ws.Range("L939").Value = "Segunda"
$ws.Range("M939").Value = 20
$ws.Range("N939").Value = 120000
$ws.Range("O939").Value = 130000
$ws.Range("P939").Value = 125000
$ws.Range("Q939").Value = "$/bins (400 kilos)"
$ws.Range("R939").Value = "Provincia de Limarí"
$ws.Range("S939").Value = 312
$ws.Range("T939").Value = 400

Write-Output ("Used range after edit: " + $ws.UsedRange.Address())
